# Auto-generated Excel COM-interop script.
# Adds a new weekly cohort (period 20210917-20210924) and a new date column
# (20220304) to every sheet, inserting a new cohort row before the SUM row and
# recomputing the SUM row for the affected columns.

$wb = $excel.ActiveWorkbook

# ---------- Sheet 1 ----------
$ws = $wb.Worksheets.Item(1)

# 1) Add the new date-header column (AM1), forcing text storage so it
#    matches the existing "yyyymmdd" text headers rather than becoming a number.
$ws.Cells.Item(1,39).NumberFormat = "@"
$ws.Cells.Item(1,39).Value2 = "20220304"
$ws.Cells.Item(1,38).Copy()
$ws.Cells.Item(1,39).PasteSpecial(-4122)

# 2) Insert a new row above the SUM row (old row 16 -> row 17) for the new cohort.
$ws.Rows.Item(16).Insert()

# 3) Style the new A16 label cell like the other period-label cells, then set its text.
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value2 = "20210917-20210924"

# 4) Fill the new cohort row (P16:AM16) with its constant weekly distribution value.
$newRowValue = 422993.89125
for ($col = 16; $col -le 39; $col++) {
    $ws.Cells.Item(16, $col).Value2 = $newRowValue
}

# 5) Recompute the SUM row (now row 17) for columns P:AM to include the new cohort.
$sumVals = @(
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    12545642.68708333,
    11605846.59791666,
    10666050.50875,
    9726254.419583332,
    8786458.330416666,
    7744163.402500001,
    6701868.474583334,
    5659573.546666668,
    4617278.61875,
    3780204.3825,
    2943130.14625,
    2106055.91,
    1268981.67375,
    845987.7825000001,
    422993.89125
)
for ($i = 0; $i -lt $sumVals.Length; $i++) {
    $ws.Cells.Item(17, 16 + $i).Value2 = $sumVals[$i]
}

# ---------- Sheet 2 ----------
$ws = $wb.Worksheets.Item(2)

# 1) Add the new date-header column (AM1), forcing text storage so it
#    matches the existing "yyyymmdd" text headers rather than becoming a number.
$ws.Cells.Item(1,39).NumberFormat = "@"
$ws.Cells.Item(1,39).Value2 = "20220304"
$ws.Cells.Item(1,38).Copy()
$ws.Cells.Item(1,39).PasteSpecial(-4122)

# 2) Insert a new row above the SUM row (old row 16 -> row 17) for the new cohort.
$ws.Rows.Item(16).Insert()

# 3) Style the new A16 label cell like the other period-label cells, then set its text.
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value2 = "20210917-20210924"

# 4) Fill the new cohort row (P16:AM16) with its constant weekly distribution value.
$newRowValue = 701978.0920833334
for ($col = 16; $col -le 39; $col++) {
    $ws.Cells.Item(16, $col).Value2 = $newRowValue
}

# 5) Recompute the SUM row (now row 17) for columns P:AM to include the new cohort.
$sumVals = @(
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    13073767.60958333,
    12131684.27625,
    11244142.60958333,
    10356600.94291667,
    9469059.27625,
    8596392.609583333,
    7723725.942916665,
    6851059.276249999,
    5978392.609583333,
    5051184.27625,
    4069434.276250001,
    3087684.276250001,
    2105934.27625,
    1403956.184166667,
    701978.0920833334
)
for ($i = 0; $i -lt $sumVals.Length; $i++) {
    $ws.Cells.Item(17, 16 + $i).Value2 = $sumVals[$i]
}

# ---------- Sheet 3 ----------
$ws = $wb.Worksheets.Item(3)

# 1) Add the new date-header column (AM1), forcing text storage so it
#    matches the existing "yyyymmdd" text headers rather than becoming a number.
$ws.Cells.Item(1,39).NumberFormat = "@"
$ws.Cells.Item(1,39).Value2 = "20220304"
$ws.Cells.Item(1,38).Copy()
$ws.Cells.Item(1,39).PasteSpecial(-4122)

# 2) Insert a new row above the SUM row (old row 16 -> row 17) for the new cohort.
$ws.Rows.Item(16).Insert()

# 3) Style the new A16 label cell like the other period-label cells, then set its text.
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value2 = "20210917-20210924"

# 4) Fill the new cohort row (P16:AM16) with its constant weekly distribution value.
$newRowValue = 1513058.30375
for ($col = 16; $col -le 39; $col++) {
    $ws.Cells.Item(16, $col).Value2 = $newRowValue
}

# 5) Recompute the SUM row (now row 17) for columns P:AM to include the new cohort.
$sumVals = @(
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    28814898.09625001,
    26791921.16416667,
    24768944.23208334,
    22745967.3,
    20722990.36791667,
    18700013.43583333,
    16677036.50375,
    14654059.57166667,
    12631082.63958333,
    10608105.7075,
    8585128.775416667,
    6562151.843333334,
    4539174.91125,
    3026116.6075,
    1513058.30375
)
for ($i = 0; $i -lt $sumVals.Length; $i++) {
    $ws.Cells.Item(17, 16 + $i).Value2 = $sumVals[$i]
}

# ---------- Sheet 4 ----------
$ws = $wb.Worksheets.Item(4)

# 1) Add the new date-header column (AM1), forcing text storage so it
#    matches the existing "yyyymmdd" text headers rather than becoming a number.
$ws.Cells.Item(1,39).NumberFormat = "@"
$ws.Cells.Item(1,39).Value2 = "20220304"
$ws.Cells.Item(1,38).Copy()
$ws.Cells.Item(1,39).PasteSpecial(-4122)

# 2) Insert a new row above the SUM row (old row 16 -> row 17) for the new cohort.
$ws.Rows.Item(16).Insert()

# 3) Style the new A16 label cell like the other period-label cells, then set its text.
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value2 = "20210917-20210924"

# 4) Fill the new cohort row (P16:AM16) with its constant weekly distribution value.
$newRowValue = 3632703.442083333
for ($col = 16; $col -le 39; $col++) {
    $ws.Cells.Item(16, $col).Value2 = $newRowValue
}

# 5) Recompute the SUM row (now row 17) for columns P:AM to include the new cohort.
$sumVals = @(
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    74729899.39625001,
    69410583.64041668,
    64091267.88458335,
    58771952.12875002,
    53452636.37291668,
    48133320.61708335,
    42814004.86125001,
    37494689.10541667,
    32175373.34958333,
    26856057.59375,
    21536741.83791666,
    16217426.08208333,
    10898110.32625,
    7265406.884166666,
    3632703.442083333
)
for ($i = 0; $i -lt $sumVals.Length; $i++) {
    $ws.Cells.Item(17, 16 + $i).Value2 = $sumVals[$i]
}

# ---------- Sheet 5 ----------
$ws = $wb.Worksheets.Item(5)

# 1) Add the new date-header column (AM1), forcing text storage so it
#    matches the existing "yyyymmdd" text headers rather than becoming a number.
$ws.Cells.Item(1,39).NumberFormat = "@"
$ws.Cells.Item(1,39).Value2 = "20220304"
$ws.Cells.Item(1,38).Copy()
$ws.Cells.Item(1,39).PasteSpecial(-4122)

# 2) Insert a new row above the SUM row (old row 16 -> row 17) for the new cohort.
$ws.Rows.Item(16).Insert()

# 3) Style the new A16 label cell like the other period-label cells, then set its text.
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value2 = "20210917-20210924"

# 4) Fill the new cohort row (P16:AM16) with its constant weekly distribution value.
$newRowValue = 2257689781.774583
for ($col = 16; $col -le 39; $col++) {
    $ws.Cells.Item(16, $col).Value2 = $newRowValue
}

# 5) Recompute the SUM row (now row 17) for columns P:AM to include the new cohort.
$sumVals = @(
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    40771010835.21375,
    37937849044.38958,
    35104687253.56541,
    32271525462.74125,
    29438363671.91708,
    26605201881.09291,
    23772040090.26875,
    20938878299.44458,
    18105716508.62042,
    15272554717.79625,
    12439392926.97208,
    9606231136.147917,
    6773069345.32375,
    4515379563.549167,
    2257689781.774583
)
for ($i = 0; $i -lt $sumVals.Length; $i++) {
    $ws.Cells.Item(17, 16 + $i).Value2 = $sumVals[$i]
}

